# Auto-generated edit script: update cryptos list values (commit: "Updated cryptos list on Mon Dec  2 22:37:39 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '95.224.60'
$ws.Range('E2').Value = '  -2.36%  '

# Row 3
$ws.Range('D3').Value = '3.606.93'
$ws.Range('E3').Value = '  -3.22%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '2.73'
$ws.Range('E4').Value = '  +25.68%  '

# Row 5
$ws.Range('E5').Value = '  +0.16%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '223.48'
$ws.Range('E6').Value = '  -6.17%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '639.83'
$ws.Range('E7').Value = '  -2.71%  '

# Row 8
$ws.Range('E8').Value = '  -5.48%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.18'
$ws.Range('E9').Value = '  +4.80%  '

# Row 10
$ws.Range('E10').Value = '  +0.01%  '

# Row 11
$ws.Range('D11').Value = '3.601.88'
$ws.Range('E11').Value = '  -3.32%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '50.43'
$ws.Range('E12').Value = '  +12.61%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.216'
$ws.Range('E13').Value = '  +4.45%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000293'
$ws.Range('E14').Value = '  -7.04%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.49'
$ws.Range('E15').Value = '  -5.27%  '

# Row 16
$ws.Range('D16').Value = '4.279.80'
$ws.Range('E16').Value = '  -3.23%  '

# Row 17
$ws.Range('D17').Value = '95.161.35'
$ws.Range('E17').Value = '  -2.18%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '24.48'
$ws.Range('E18').Value = '  +29.69%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.11'
$ws.Range('E19').Value = '  +2.75%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.72'
$ws.Range('E20').Value = '  +4.67%  '

# Row 21
$ws.Range('D21').Value = '3.607.63'
$ws.Range('E21').Value = '  -2.97%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.291'
$ws.Range('E22').Value = '  +36.93%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.532'
$ws.Range('E23').Value = '  -1.35%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '136.11'
$ws.Range('E24').Value = '  +16.20%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '530.22'
$ws.Range('E25').Value = '  +0.26%  '

# Row 26
$ws.Range('E26').Value = '  -5.92%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.01'
$ws.Range('E27').Value = '  +1.59%  '

# Row 28
$ws.Range('E28').Value = '  -9.94%  '

# Row 29
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '3.790.36'
$ws.Range('E29').Value = '  -3.44%  '

# Row 30
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '13.14'
$ws.Range('E30').Value = '  -2.39%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.30'
$ws.Range('E31').Value = '  +3.93%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.14'
$ws.Range('E32').Value = '  +3.26%  '

# Row 33
$ws.Range('E33').Value = '  +0.10%  '

# Row 34
$ws.Range('E34').Value = '  +2.07%  '

# Row 35
$ws.Range('B35').Value = 'PolygonEcosystemToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.635'
$ws.Range('E35').Value = '  +6.05%  '

# Row 36
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '33.61'
$ws.Range('E36').Value = '  +1.43%  '

# Row 37
$ws.Range('E37').Value = '  -4.05%  '

# Row 38
$ws.Range('E38').Value = '  +0.21%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0556'
$ws.Range('E39').Value = '  +21.43%  '

# Row 40
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.03%  '

# Row 41
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.56'
$ws.Range('E41').Value = '  -2.34%  '

# Row 42
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.28'
$ws.Range('E42').Value = '  +6.48%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '588.71'
$ws.Range('E43').Value = '  -8.07%  '

# Row 44
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.502'
$ws.Range('E44').Value = '  +0.59%  '

# Row 45
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.01'
$ws.Range('E45').Value = '  +4.80%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.93'
$ws.Range('E46').Value = '  +0.27%  '

# Row 47
$ws.Range('B47').Value = 'ImmutableX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.00'
$ws.Range('E47').Value = '  -0.78%  '

# Row 48
$ws.Range('B48').Value = 'Kaspa'
$ws.Range('C48').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.156'
$ws.Range('E48').Value = '  -7.42%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.26'
$ws.Range('E49').Value = '  +5.39%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '233.56'
$ws.Range('E50').Value = '  +11.88%  '

# Row 51
$ws.Range('E51').Value = '  -2.40%  '
